$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(45973, 8901.46666004607, 8427.49659142224, 13084, 4801.00979475575, 6.02109942408318)
    3  = @(45974, 8841.03975577193, 8255.48077901061, 9004, 4938.78783890964, 174.594525746677)
    4  = @(45975, 8729.69078950592, 7303.41485018583, 9004, 4848.23652232565, 131.152140521312)
    5  = @(45976, 3069.61416164093, 4833.23816705591, 9004, 4370.62946535144, 8.32781801697259)
    6  = @(45977, 2814.89706303287, 4533.01905785672, 9004, 4220.24733008045, -10.4472338359512)
    7  = @(45978, 10174.4835962506, 8084.29583335618, 9004, 5864.29314492894, 206.02454076188)
    8  = @(45979, 10116.508538002, 8529.8290446245, 9004, 5829.53432565045, 223.140140428123)
    9  = @(45980, 10116.508538002, 9156.452994414, 9004, 5829.53432565045, 249.249471669352)
    10 = @(45981, 10116.508538002, 10147.7303033335, 9004, 5829.53432565045, 290.552692874332)
    11 = @(45982, 10116.508538002, 9519.66233446131, 9004, 5829.53432565045, 264.383194171323)
    12 = @(45983, 4000.30754528634, 7171.32852928911, 9004, 5477.69526251754, 151.875991325277)
    13 = @(45984, 3853.28814954241, 7135.32181685114, 9004, 5469.69176310149, 150.042232498026)
    14 = @(45985, 10980.8696790314, 10491.9520115779, 9004, 6409.17212077949, 329.046838848223)
    15 = @(45986, 10980.8696790314, 10242.0605086479, 9004, 6409.17212077949, 318.634692892808)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
